# Regenerate the order workbook with updated distance/size codes.
#
# The original experiment used Distance codes D51/D64/D80 and a Size code
# S30; this revision renumbers them to D55/D69/D86 and S31 respectively
# (S20/S25 are unchanged). The remap is applied as a literal, case-sensitive
# text substitution across every textual cell in the sheet (condition
# labels, stimulus filenames, and the standalone Distance/Size lookup
# columns all embed these tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

# Pull the whole used range into memory as a single 2-D array so the
# substitution runs as one read + one write instead of one COM round-trip
# per cell.
$arr = $used.Value2

# Order matters only in that none of the replacement tokens themselves
# contain another pattern's search token, so a single left-to-right pass
# per cell is safe and idempotent.
$map = [ordered]@{
    'D51' = 'D55'
    'D80' = 'D86'
    'D64' = 'D69'
    'S30' = 'S31'
}

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $v = $arr[$r, $c]
        if ($v -is [string]) {
            $nv = $v
            foreach ($key in $map.Keys) {
                $nv = $nv -creplace $key, $map[$key]
            }
            if ($nv -cne $v) {
                $arr[$r, $c] = $nv
            }
        }
    }
}

$used.Value2 = $arr
